$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from 3 to 2
$ws.Range("B2").Value = 2

# Add new row 3: A3 = 0 (styled like A2/B1), B3 = 1
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1

# Copy style (font, border, alignment) from A2 to A3
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
